$d = $word.ActiveDocument

# 1) Change the Varchar length for LYRICS column from 4000 to 8000
$d.Content.Find.Execute("4000", $false, $false, $false, $false, $false, $true, 1, $false, "8000", 2)

# 2) Change the cached TIME field result date from 2011-02-22 to 2011-02-23
$d.Content.Find.Execute("2011-02-22", $false, $false, $false, $false, $false, $true, 1, $false, "2011-02-23", 2)
